# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) and on each
#   per-language sheet's Status column.
# - Shrink the now-narrower Status column(s) to match the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 12.5   # -> stored column width ~13.33 chars (was ~17.22)

# --- Overview sheet: "zh-cn" (col E) and "de-de" (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: "Status" column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: "Status" column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
